$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(6)
$shape.Left = 26.784961629921263
$shape.Top = 90.83889463779526
